$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shp = $s.Shapes.Item(3)

$tr = $shp.TextFrame.TextRange
$beforeText = $tr.Text
$beforeLen = $beforeText.Length

# Fix spcAft of the "Once we have..." paragraph (1000 -> 0)
$onceIdx = $beforeText.IndexOf("Once we have the data and model")
$onceRange = $tr.Characters($onceIdx + 1, 10)
$onceRange.ParagraphFormat.SpaceAfter = 0

# Append a new paragraph with the FAIR data practices bullet
$newPara = $tr.InsertAfter("`rFAIR (Findable, Accessible, Interoperable, and Reusable) data practices.")

$fairStart = $beforeLen + 2
$fairRun = $tr.Characters($fairStart, 4)
$fairRun.Font.Bold = $true
$fairRun.Font.Size = 16

$restRun = $tr.Characters($fairStart + 4, 70)
$restRun.Font.Size = 16

# New paragraph should also have spcAft = 0
$newParaRange = $tr.Characters($fairStart, 4)
$newParaRange.ParagraphFormat.SpaceAfter = 0

# Resize the shape (a:ext cy 1906500 -> 2335800 EMU = 150.1181pt -> 183.9213pt)
# Must be set after the text edits since the shape has spAutoFit and would
# otherwise get recalculated based on the new text content.
$shp.Height = 2335800 / 12700
